$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144:176 down to 145:177
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new record
$ws.Range("A144").Value = 9
$ws.Range("B144").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C144").Value = "Metropolitana"
$ws.Range("D144").Value = 45135
$ws.Range("E144").Value = 13
$ws.Range("F144").Value = 100112022
$ws.Range("G144").Value = "Arveja Verde"
$ws.Range("H144").Value = "Perfection"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 52
$ws.Range("K144").Value = 18000
$ws.Range("L144").Value = 19000
$ws.Range("M144").Value = 18500
$ws.Range("N144").Value = "`$/malla 25 kilos"
$ws.Range("O144").Value = "Provincia de Limarí"
$ws.Range("P144").Value = 740
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
